$d = $word.ActiveDocument

# --- Paragraph 1: "The respondent may upload..." -> "The defendant may upload..."
$rng1 = $d.Content.Duplicate
$rng1.Find.Execute("respondent", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng1.Find.Found) {
  $rng1.Text = "defendant"
  # Toggle a character property and flip it back off again. This makes Word
  # keep the replaced word as its own run (instead of silently re-merging it
  # back into the identically-formatted runs on either side), while leaving
  # the run's final formatting identical to its neighbours, matching the
  # run split produced by editing just this one word in the sentence.
  $rng1.Bold = 1
  $rng1.Bold = $false
}

# --- Paragraph 2: "The applicant may upload..." -> "The claimant may upload..."
$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("applicant", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
  $rng2.Text = "claimant"
  $rng2.Bold = 1
  $rng2.Bold = $false
}
